$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two target-bandwidth input cells; all dependent formulas
# (K4, K7-K11, K15, K17-K22, K28, K29-K33, O28, O32, O33, B34-B37)
# recalculate automatically.
$ws.Range("K3").Value = 400000
$ws.Range("K14").Value = 300000

# Selection / view changes captured in the sheetView.
$ws.Range("B34:B37").Select()
$excel.ActiveWindow.ScrollRow = 7

$ws.Range("K9:K10").WrapText = $false
